$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.475.86"
$ws.Cells.Item(2, 5).Value = "  +0.58%  "
$ws.Cells.Item(3, 4).Value = "2.104.94"
$ws.Cells.Item(3, 5).Value = "  +1.12%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.007"
$ws.Cells.Item(4, 5).Value = "  +0.65%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "333.79"
$ws.Cells.Item(5, 5).Value = "  +1.63%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.006"
$ws.Cells.Item(6, 5).Value = "  +0.60%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.5220"
$ws.Cells.Item(7, 5).Value = "  -0.01%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.4524"
$ws.Cells.Item(8, 5).Value = "  +4.79%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "53.41"
$ws.Cells.Item(9, 5).Value = "  +14.29%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08915"
$ws.Cells.Item(10, 5).Value = "  +0.98%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.183"
$ws.Cells.Item(11, 5).Value = "  +1.82%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "24.08"
$ws.Cells.Item(12, 5).Value = "  -1.50%  "
$ws.Cells.Item(13, 4).Value = "2.099.83"
$ws.Cells.Item(13, 5).Value = "  +0.56%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.816"
$ws.Cells.Item(14, 5).Value = "  +1.41%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "8.025"
$ws.Cells.Item(15, 5).Value = "  +4.70%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "96.75"
$ws.Cells.Item(16, 5).Value = "  +1.36%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001142"
$ws.Cells.Item(17, 5).Value = "  +1.66%  "
$ws.Cells.Item(18, 5).Value = "  +0.50%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06653"
$ws.Cells.Item(19, 5).Value = "  +0.30%  "
$ws.Cells.Item(20, 5).Value = "  +1.84%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "1.005"
$ws.Cells.Item(21, 5).Value = "  +0.46%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.335"
$ws.Cells.Item(22, 5).Value = "  +0.47%  "
$ws.Cells.Item(23, 4).Value = "30.527.10"
$ws.Cells.Item(23, 5).Value = "  +0.57%  "
$ws.Cells.Item(24, 5).Value = "  +1.07%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.354"
$ws.Cells.Item(25, 5).Value = "  +2.37%  "
$ws.Cells.Item(26, 4).Value = "2.350.15"
$ws.Cells.Item(26, 5).Value = "  +0.76%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "22.29"
$ws.Cells.Item(27, 5).Value = "  -0.39%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "162.72"
$ws.Cells.Item(28, 5).Value = "  +0.44%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.524"
$ws.Cells.Item(29, 5).Value = "  -2.71%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "133.25"
$ws.Cells.Item(30, 5).Value = "  +1.45%  "
$ws.Cells.Item(31, 5).Value = "  +1.35%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.1070"
$ws.Cells.Item(32, 5).Value = "  +0.04%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.657"
$ws.Cells.Item(33, 5).Value = "  +0.41%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.436"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.947"
$ws.Cells.Item(35, 5).Value = "  +2.24%  "
$ws.Cells.Item(36, 5).Value = "  +4.98%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.796"
$ws.Cells.Item(37, 5).Value = "  +6.34%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02590"
$ws.Cells.Item(38, 5).Value = "  +0.76%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.06845"
$ws.Cells.Item(39, 5).Value = "  +2.44%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.2296"
$ws.Cells.Item(40, 5).Value = "  +1.49%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "12.71"
$ws.Cells.Item(41, 5).Value = "  +0.20%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.6870"
$ws.Cells.Item(42, 5).Value = "  +0.80%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.248"
$ws.Cells.Item(43, 5).Value = "  +0.22%  "
$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "14.12"
$ws.Cells.Item(44, 5).Value = "  +0.55%  "
$ws.Cells.Item(45, 2).Value = "NEARProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.311"
$ws.Cells.Item(45, 5).Value = "  +4.95%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.6357"
$ws.Cells.Item(46, 5).Value = "  -0.14%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "3.659"
$ws.Cells.Item(47, 5).Value = "  +1.31%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.00000000353"
$ws.Cells.Item(48, 5).Value = "  +22.44%  "
$ws.Cells.Item(49, 5).Value = "  -0.31%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "83.53"
$ws.Cells.Item(50, 5).Value = "  +2.42%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.207"
$ws.Cells.Item(51, 5).Value = "  +1.80%  "